$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("C7").Value = 100
$ws.Range("D8").Value = 0.063

# Add new row for Muscles
$ws.Range("A13").Value = "Muscles"
$ws.Range("B13").Value = 40
$ws.Range("C13").Value = 27.8
$ws.Range("D13").Value = 0.0143
$ws.Range("E13").Value = 0.02

# Update selection to match target state
$ws.Range("C8").Select()
